# Update the cryptocurrency prices / 1h volume % figures on the active sheet
# to reflect the latest scrape, per the GitHub Actions commit.
#
# NOTE: several Price (column D) values are plain decimal numbers
# (e.g. "673.61"). Excel's COM layer auto-converts a Range.Value string that
# parses as a number into a numeric cell. The source data models these as
# text (t="inlineStr" in the original file), so for those cells we force the
# NumberFormat to Text ("@") before assigning the value to keep them as
# strings, matching values like "69.654.05" that already stay text because
# they aren't valid numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.654.05"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").Value = "3.705.87"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "673.61"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.05"
$ws.Range("E6").Value = "  +2.49%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  +1.11%  "

$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("E10").Value = "  +1.34%  "

$ws.Range("E11").Value = "  +1.68%  "

$ws.Range("E12").Value = "  +1.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.97"
$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("D14").Value = "3.713.24"
$ws.Range("E14").Value = "  +1.63%  "

$ws.Range("D15").Value = "69.636.40"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("E16").Value = "  +1.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.37"
$ws.Range("E17").Value = "  +2.66%  "

$ws.Range("E18").Value = "  +2.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "474.51"
$ws.Range("E19").Value = "  +0.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.82"
$ws.Range("E20").Value = "  -2.38%  "

$ws.Range("E21").Value = "  +0.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "80.41"
$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("D23").Value = "3.854.00"
$ws.Range("E23").Value = "  +0.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000128"
$ws.Range("E24").Value = "  +5.23%  "

$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("E26").Value = "  +1.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("E27").Value = "  +0.59%  "

$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.74"
$ws.Range("E29").Value = "  -0.75%  "

# Rows 30/31 swap their coin data: Kaspa moves up to row 30, ImmutableX down to row 31.
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.169"
$ws.Range("E30").Value = "  +7.89%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.01"
$ws.Range("E31").Value = "  +1.41%  "

$ws.Range("E32").Value = "  +0.32%  "

$ws.Range("E33").Value = "  -0.36%  "

$ws.Range("E34").Value = "  +0.38%  "

$ws.Range("D35").Value = "3.696.00"
$ws.Range("E35").Value = "  +0.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.59"
$ws.Range("E36").Value = "  +4.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.13"
$ws.Range("E37").Value = "  +1.15%  "

$ws.Range("E39").Value = "  +1.04%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  +1.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "174.01"
$ws.Range("E42").Value = "  +2.84%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.09"
$ws.Range("E44").Value = "  -0.95%  "

$ws.Range("E45").Value = "  +2.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000280"
$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("E47").Value = "  +2.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.88"
$ws.Range("E48").Value = "  +3.63%  "

$ws.Range("E49").Value = "  -1.22%  "

$ws.Range("E50").Value = "  +1.55%  "

$ws.Range("E51").Value = "  +0.55%  "
